# "add week 13 14, logreg"
#
# Updates the schedule for week 12/13 (rows 13-14: Logistic Regression /
# Multiple Regression material gets split out and the reading list is
# expanded) and reshuffles the "In class presentations" note from
# row 17 (G17) / Final Posters note (H17) along with the Special Analysis
# Topics (D16) and Poster Revisions (D17) text, plus moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (week 12: "Categorical Predictors" / multiple+logistic regression notes) ---
$ws.Range("D13").Value = "Categorical Predictors`n* Log and logit transformations "
$ws.Range("G13").Value = "Multiple Regression Notes [[HTML]](docs/lec07_MLR.html) [[PDF]](docs/lec07_MLR.pdf)`n* Logistic Regression Notes [[HTML]](docs/lec08_LogReg.html) [[PDF]](docs/lec08_LogReg.pdf)`n* PMA5 9.3 (dummy variables)`n* PMA5 Ch 6.9 (transformations)`n* PMA5 Ch 12 (Logistic Regression)`n"
$ws.Rows(13).RowHeight = 173.25

# --- Row 14 (week 13: "Model building/ fit" / variable selection reading) ---
$ws.Range("D14").Value = "Model building/ fit"
$ws.Range("G14").Value = "PMA5 Ch 8 (Variable Selection)"
$ws.Rows(14).RowHeight = 31.5

# --- Row 16 (week 14: Special Analysis Topics) ---
$ws.Range("D16").Value = "Special Analysis Topics`n* Poster design"

# --- Row 17 (week 15: Poster Revisions / in class presentations / final posters) ---
$ws.Range("D17").Value = "Poster Revisions`n* Poster Presentations"
$ws.Range("G17").Value = "In class presentations on 12/7. "
$ws.Range("H17").Value = "Final Posters (as printed)  Due 12/7 EOD"

# --- Move active selection to G15 ---
$ws.Range("G15").Select()
